# CSAT_Performance_Reports.xlsx update
# - Daywise_Report: update MTD totals, append two new daily rows (10/05, 10/06)
# - Agentwise_Report: insert two new agents (SBM990, VPS193) keeping alpha order
# - Daywise_Agent_Performance: append two matching daywise/agent rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Daywise_Report  (Table_Daywise_Report, A1:F10 -> A1:F12)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Daywise_Report")
$lo1 = $ws1.ListObjects.Item("Table_Daywise_Report")

# Update the MTD summary row (row 2)
$ws1.Range("D2").Value = 11
$ws1.Range("E2").Value = 21
$ws1.Range("F2").Value = 2.14

# Grow the table by two rows (appended at the end)
$lo1.ListRows.Add() | Out-Null
$lo1.ListRows.Add() | Out-Null

# Row 3 (10/01) has the exact formatting template we need (date style,
# red CSAT3 fill, green Score fill) for both brand new rows.
$ws1.Range("A3:F3").Copy()
$ws1.Range("A11:F11").PasteSpecial(-4122)
$ws1.Range("A3:F3").Copy()
$ws1.Range("A12:F12").PasteSpecial(-4122)

$ws1.Range("A11").Value = 45940
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 0
$ws1.Range("D11").Value = 1
$ws1.Range("E11").Value = 1
$ws1.Range("F11").Value = 3

$ws1.Range("A12").Value = 45941
$ws1.Range("B12").Value = 0
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = 1
$ws1.Range("E12").Value = 1
$ws1.Range("F12").Value = 3

# ---------------------------------------------------------------------------
# Sheet 2: Agentwise_Report  (Table_Agentwise_Report, A1:F18 -> A1:F20)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Agentwise_Report")
$lo2 = $ws2.ListObjects.Item("Table_Agentwise_Report")

# Insert "SBM990" right before the current row 11 ("VPS111"), keeping the
# agent list alphabetically sorted.
$ws2.Rows.Item(11).Insert()
$lo2.Resize($ws2.Range("A1:F19"))

$ws2.Range("A11").Value = "SBM990"
$ws2.Range("B11").Value = 0
$ws2.Range("C11").Value = 0
$ws2.Range("D11").Value = 1
$ws2.Range("E11").Value = 1
$ws2.Range("F11").Value = 3

# Insert "VPS193" right before the now-shifted row 13 ("VPS214") -- VPS111
# sits at row 12 after the first insert above.
$ws2.Rows.Item(13).Insert()
$lo2.Resize($ws2.Range("A1:F20"))

$ws2.Range("A13").Value = "VPS193"
$ws2.Range("B13").Value = 0
$ws2.Range("C13").Value = 0
$ws2.Range("D13").Value = 1
$ws2.Range("E13").Value = 1
$ws2.Range("F13").Value = 3

# ---------------------------------------------------------------------------
# Sheet 4: Daywise_Agent_Performance  (Table_Daywise_Agent_Performance, A1:G20 -> A1:G22)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Daywise_Agent_Performance")
$lo4 = $ws4.ListObjects.Item("Table_Daywise_Agent_Performance")

$lo4.ListRows.Add() | Out-Null
$lo4.ListRows.Add() | Out-Null

# Row 13 has the exact formatting template needed (date style, plain
# CSAT cells, green Score fill) for both brand new rows.
$ws4.Range("A13:G13").Copy()
$ws4.Range("A21:G21").PasteSpecial(-4122)
$ws4.Range("A13:G13").Copy()
$ws4.Range("A22:G22").PasteSpecial(-4122)

$ws4.Range("A21").Value = 45940
$ws4.Range("B21").Value = "SBM990"
$ws4.Range("C21").Value = 0
$ws4.Range("D21").Value = 0
$ws4.Range("E21").Value = 1
$ws4.Range("F21").Value = 1
$ws4.Range("G21").Value = 3

$ws4.Range("A22").Value = 45941
$ws4.Range("B22").Value = "VPS193"
$ws4.Range("C22").Value = 0
$ws4.Range("D22").Value = 0
$ws4.Range("E22").Value = 1
$ws4.Range("F22").Value = 1
$ws4.Range("G22").Value = 3
